$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The item list (A/B/H/L/N columns, rows 4-23) is kept alphabetically sorted
# by product name (column B). Two new products were added:
#   - ISOMOISTY NASAL SPRAY 30 ML   (sorts between HIBIOTIC and MEGALASE)
#   - TELFAST 30MG/5ML SUSP. 100 ML (sorts between PROXIMOL and VOLTAREN)
# Inserting them re-flows every row from the insertion point onward, and the
# table grows from 20 data rows (rows 4-23) to 22 data rows (rows 4-25).
# The previous "totals" row (24) and "footer" row (25) therefore move down
# to rows 26 and 27.
# ---------------------------------------------------------------------------

# Step 1: move the totals ("K" sum) row from 24 -> 26, and the footer row
# from 25 -> 27, preserving their formatting. Do this first (before we
# overwrite rows 24/25 with new product data) and before we unmerge, so the
# formats are still intact to copy from.
$ws.Range("A24:N24").Copy()
$ws.Range("A26:N26").PasteSpecial(-4122)

$ws.Range("A25:N25").Copy()
$ws.Range("A27:N27").PasteSpecial(-4122)

# Remove the old merges tied to rows 24/25 (now stale) ...
$ws.Range("K24:N24").UnMerge()
$ws.Range("A25:E25").UnMerge()
$ws.Range("F25:G25").UnMerge()
$ws.Range("I25:N25").UnMerge()

# ... and recreate them on the new rows 26/27.
$ws.Range("K26:N26").Merge()
$ws.Range("A27:E27").Merge()
$ws.Range("F27:G27").Merge()
$ws.Range("I27:N27").Merge()

# New totals value (old 1470.95 + 90 (ISOMOISTY) + 50 (TELFAST) = 1610.95)
$ws.Range("K26").Value = 1610.9500000000001

# Footer row content is unchanged, just relocated.
$ws.Range("A27").Value = "Tuesday, 13 January, 2026 2:35 PM"
$ws.Range("F27").Value = "1/1"
$ws.Range("I27").Value = "developed by : Abdelaziz Talaat"

# Step 2: create the two brand-new data rows (24 and 25) by cloning the
# formatting of the last existing data row (23), then merging like every
# other data row.
$ws.Range("A23:N23").Copy()
$ws.Range("A24:N24").PasteSpecial(-4122)
$ws.Range("A23:N23").Copy()
$ws.Range("A25:N25").PasteSpecial(-4122)

$ws.Range("B24:G24").Merge()
$ws.Range("H24:K24").Merge()
$ws.Range("L24:M24").Merge()
$ws.Range("B25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()

# Step 3: rewrite the data rows 4-25 with the new, alphabetically sorted
# product list. Rows 4-13 are unaffected by the insert, rows 14-25 absorb
# the shift.
$names = @(
  "ANTINAL 220MG/5ML 60ML SUSP.",
  "AVIVAVASC 5/160MG 28 F.C. TAB.",
  "CARBAMIDE 10% CREAM 30 GM",
  "CONTAFEVER N 200MG/5ML SUSP. 120ML",
  "DEPOVIT B12-1000MCG/ML 5 I.M. AMP",
  "DOLO-D PLUS ORAL SUSP. 115 ML",
  "EXOSIRYLIC 500 MG 20 F.C.TABS.",
  "FORBUDES 400/12MCG 60 INHALATION CAPS.+INHALER",
  "GLUCOVANCE 500/5MG 30 F.C.TAB.",
  "HIBIOTIC N 600MG SUSP. 80 ML",
  "ISOMOISTY NASAL SPRAY 30 ML",
  "MEGALASE SYRUP 125 ML",
  "MINALAX 10 TABLETS",
  "NOSTAMINE EYE/NOSE DROPS 15 ML",
  "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML",
  "OTRIVIN 0.1% ADULT NASAL DROPS 15 ML",
  "PROXIMOL 0.4MG 40 TAB",
  "TELFAST 30MG/5ML SUSP. 100 ML",
  "VOLTAREN 75MG/3ML 3 AMP.",
  "WATER FOR INJECTION AMP. 5 ML",
  "سرنجات 3 سم",
  "معجون سيجنال 25 مل"
)
$balances = @("4:0","0:0","2:0","8:0","2:2","3:0","0:1","1:0","0:0","1:0","0:0","2:0","7:0","9:0","0:0","4:0","0:1","3:0","5:2","7762:0","-2:0","2:0")
$prices = @(24,120.28,40,66,85,41,194,55.67,74,92,90,31,36,440,48,24,34,50,17,5,24,20)
$counts = @("1:0","1:0","1:0","2:0","1:0","1:0","1:0","0:0","1:0","1:0","1:0","1:0","2:0","20:0","2:0","1:0","1:0","1:0","0:0","2:0","12:0","1:0")

for ($i = 0; $i -lt $names.Count; $i++) {
  $r = 4 + $i
  $ws.Range("A$r").Value = $i + 1
  $ws.Range("B$r").Value = $names[$i]
  $ws.Range("H$r").Value = $balances[$i]
  $ws.Range("L$r").Value = $prices[$i]
  $ws.Range("N$r").Value = $counts[$i]
}
